# Automatische test-sync: 2025-08-14 22:13:50
# Adds a new log entry (row 41) to the "Logs" sheet, extends the
# conditional formatting ranges to cover the new row, and updates the
# "Dashboard" sheet's aggregate count for the affected category.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new row -----------------------------------
$ws = $wb.Worksheets.Item("Logs")

$newRow = 41
$ws.Cells.Item($newRow, 1).Value = "Opvolging retour"
$ws.Cells.Item($newRow, 2).Value = "kwaliteit@testbedrijf123.nl"
$ws.Cells.Item($newRow, 3).Value = "Hebben jullie al nieuws over mijn retour?"
$ws.Cells.Item($newRow, 4).Value = "Intern verzoek / Actie voor medewerker"
$ws.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar klantenservice@testbedrijf123.nl."
$ws.Cells.Item($newRow, 6).Value = "2025-08-14 22:13:04"
$ws.Cells.Item($newRow, 7).Value = "Nee"
$ws.Cells.Item($newRow, 8).Value = "Ja"
$ws.Cells.Item($newRow, 9).Value = "Nee"
$ws.Cells.Item($newRow, 10).Value = "Nee"

# --- Extend the conditional formatting ranges to include the new row --
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $ws.Range($col + "2:" + $col + "40")
    $newRange = $ws.Range($col + "2:" + $col + "41")
    $cfs = $oldRange.FormatConditions
    for ($i = 1; $i -le $cfs.Count(); $i++) {
        $cfs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet: bump the count for the affected category --------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 33
